# ---------------------------------------------------------------------------
# "added data to the sample data" -- extends the "Recommender Math" workbook:
#   * Template sheet: a few matrix cells change (new purchase/view overlap),
#     and two new "Cross-action" lookup tables are added below the existing
#     matrices (rows 61-64 and 66-70 / 72-76).
#   * Action Logs sheet: a few rows' item is corrected, and a new action-log
#     row is appended.
#   * Template becomes the active sheet/tab again.
# ---------------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$tpl = $wb.Worksheets.Item("Template")
$log = $wb.Worksheets.Item("Action Logs")

# ---------------------------------------------------------------------------
# 1) Template sheet: corrected overlap cells in the existing matrices
# ---------------------------------------------------------------------------

$tpl.Range("E38").Value = 1

$tpl.Range("E44").Value = 1

$tpl.Range("E53").Value = 1

$tpl.Range("G56").Value = 1
$tpl.Range("P56").Value = 2

$tpl.Range("G57").Value = 1
$tpl.Range("P57").Value = 1

$tpl.Range("C58").Value = 2
$tpl.Range("D58").Value = 1
$tpl.Range("E58").Value = 1
$tpl.Range("F58").Value = 2
$tpl.Range("G58").Value = 2
$tpl.Range("P58").Value = 1

$tpl.Range("G59").Value = 1
$tpl.Range("P59").Value = 2

$tpl.Range("N60").Value = 1
$tpl.Range("O60").Value = 1
$tpl.Range("P60").Value = 2
$tpl.Range("Q60").Value = 1

# ---------------------------------------------------------------------------
# 2) Template sheet: new "Cross-action Recommendations" mini table (61-64)
# ---------------------------------------------------------------------------

$tpl.Range("B61").Copy()
$tpl.Range("B61:B64").PasteSpecial(-4122)   # xlPasteFormats (label style)

$tpl.Range("A61").Value = "Cross-action Recommendations"
$tpl.Range("B61").Value = "u1"
$tpl.Range("C61").Value = "galaxy, nexus"

$tpl.Range("B62").Value = "u2"
$tpl.Range("C62").Value = "iphone, ipad"

$tpl.Range("B63").Value = "u3"
$tpl.Range("C63").Value = "iphone, galaxy, ipad, nexus"

$tpl.Range("B64").Value = "u4"
$tpl.Range("C64").Value = "ipad, nexus"

$tpl.Range("A61:A62").Merge()

# ---------------------------------------------------------------------------
# 3) Template sheet: new "Cross-action Similar Items" mini tables
#    (columns of [B'A]: rows 66-70, rows of [B'A]: rows 72-76)
# ---------------------------------------------------------------------------

$tpl.Range("B66").Copy()
$tpl.Range("B66:B70").PasteSpecial(-4122)
$tpl.Range("B72:B76").PasteSpecial(-4122)

$tpl.Range("A66").Value = "Cross-action Similari Items"
$tpl.Range("B66").Value = "iphone"
$tpl.Range("C66").Value = "galaxy, ipad, nexus"

$tpl.Range("B67").Value = "ipad"
$tpl.Range("C67").Value = "iphone, galaxy, nexus"

$tpl.Range("A68").Value = "columns of [B'A]"
$tpl.Range("B68").Value = "nexus"
$tpl.Range("C68").Value = "iphone, galaxy, ipad"

$tpl.Range("B69").Value = "galaxy"
$tpl.Range("C69").Value = "iphone, ipad, nexus"

$tpl.Range("B70").Value = "surface"

$tpl.Range("A72").Value = "Cross-action Similari Items"
$tpl.Range("B72").Value = "iphone"
$tpl.Range("C72").Value = "ipad, nexus, galaxy"

$tpl.Range("B73").Value = "ipad"
$tpl.Range("C73").Value = "iphone, nexus, galaxy"

$tpl.Range("A74").Value = "rows of [B'A]"
$tpl.Range("B74").Value = "nexus"
$tpl.Range("C74").Value = "iphone, ipad, galaxy"

$tpl.Range("B75").Value = "galaxy"
$tpl.Range("C75").Value = "iphone, ipad, nexus"

$tpl.Range("B76").Value = "surface"
$tpl.Range("C76").Value = "nexus"

$tpl.Range("A66:A67").Merge()
$tpl.Range("A72:A73").Merge()

# Wrap the long section titles in column A of the new tables
$tpl.Range("A61:A62,A66:A67,A72:A73").WrapText = $true

# ---------------------------------------------------------------------------
# 4) Action Logs sheet: fix a few item labels, append one new log row
# ---------------------------------------------------------------------------

$log.Range("C3").Value = "nexus"
$log.Range("C10").Value = "nexus"
$log.Range("C14").Value = "nexus"

$log.Range("A17").Value = "u3"
$log.Range("C17").Value = "nexus"

$log.Range("C18").Value = "iphone"

$log.Range("C19").Value = "ipad"

$log.Range("A20").Value = "u4"
$log.Range("B20").Value = "view"
$log.Range("C20").Value = "galaxy"

# ---------------------------------------------------------------------------
# 5) View state: Template becomes the active/selected sheet+tab again,
#    Action Logs' selection collapses back to the full used range.
# ---------------------------------------------------------------------------

$log.Range("A1:C20").Select()

$tpl.Activate()
$tpl.Range("C77").Select()

Write-Host "done"
